# Auto-generated edit script applying the diff to Marilith_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 500
$ws.Range("J12").Value = 500
$ws.Range("L12").Value = 500
$ws.Range("N12").Value = -840
$ws.Range("H17").Value = 1469.7
$ws.Range("J17").Value = 1469.7
$ws.Range("L17").Value = 4409.1
$ws.Range("N17").Value = -4745.1
$ws.Range("H39").Value = 650.0909
$ws.Range("I39").Value = 558.6667
$ws.Range("J39").Value = 759.8
$ws.Range("K39").Value = 1676.0001
$ws.Range("L39").Value = 2279.4
$ws.Range("M39").Value = -1380.0001
$ws.Range("N39").Value = -2871.4
$ws.Range("H51").Value = 8000
$ws.Range("I51").Value = 8000
$ws.Range("K51").Value = 8000
$ws.Range("M51").Value = -7516
$ws.Range("H98").Value = 1886.875
$ws.Range("I98").Value = 1975.8334
$ws.Range("J98").Value = 1620
$ws.Range("K98").Value = 1975.8334
$ws.Range("L98").Value = 1620
$ws.Range("M98").Value = -477.8334
$ws.Range("N98").Value = -4616
$ws.Range("H122").Value = 1886.875
$ws.Range("I122").Value = 1975.8334
$ws.Range("J122").Value = 1620
$ws.Range("K122").Value = 5927.5002
$ws.Range("L122").Value = 4860
$ws.Range("M122").Value = -3477.5002
$ws.Range("N122").Value = -9760
$ws.Range("H123").Value = 152000
$ws.Range("J123").Value = 152000
$ws.Range("L123").Value = 152000
$ws.Range("N123").Value = -161800
$ws.Range("H127").Value = 1723.375
$ws.Range("I127").Value = 1578.4
$ws.Range("J127").Value = 1965
$ws.Range("K127").Value = 4735.200000000001
$ws.Range("L127").Value = 5895
$ws.Range("M127").Value = 224.7999999999993
$ws.Range("N127").Value = -15815
$ws.Range("H132").Value = 7326
$ws.Range("I132").Value = 9501.333000000001
$ws.Range("J132").Value = 800
$ws.Range("K132").Value = 28503.999
$ws.Range("L132").Value = 2400
$ws.Range("M132").Value = -25973.999
$ws.Range("N132").Value = -7460
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 12202
$ws.Range("I14").Value = 17003.334
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 17003.334
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = -16828.334
$ws.Range("N14").Value = -5350
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 8443.25
$ws.Range("I31").Value = 6136.5
$ws.Range("J31").Value = 10750
$ws.Range("K31").Value = 6136.5
$ws.Range("L31").Value = 10750
$ws.Range("M31").Value = -5884.5
$ws.Range("N31").Value = -11254
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H76").Value = 10000
$ws.Range("I76").Value = 10000
$ws.Range("K76").Value = 10000
$ws.Range("M76").Value = -9685
$ws.Range("H79").Value = 10000
$ws.Range("I79").Value = 10000
$ws.Range("K79").Value = 10000
$ws.Range("M79").Value = -8908
$ws.Range("H107").Value = 615.1875
$ws.Range("I107").Value = 567.3333
$ws.Range("J107").Value = 758.75
$ws.Range("K107").Value = 567.3333
$ws.Range("L107").Value = 758.75
$ws.Range("M107").Value = 1352.6667
$ws.Range("N107").Value = -4598.75
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H132").Value = 3251.5833
$ws.Range("I132").Value = 2319.8333
$ws.Range("J132").Value = 4183.3335
$ws.Range("K132").Value = 6959.499899999999
$ws.Range("L132").Value = 12550.0005
$ws.Range("M132").Value = -4429.499899999999
$ws.Range("N132").Value = -17610.0005
$ws.Range("H134").Value = 1673
$ws.Range("I134").Value = 1582.2
$ws.Range("K134").Value = 4746.6
$ws.Range("M134").Value = -2211.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 208.625
$ws.Range("I86").Value = 201.8
$ws.Range("J86").Value = 220
$ws.Range("K86").Value = 605.4000000000001
$ws.Range("L86").Value = 660
$ws.Range("M86").Value = 580.5999999999999
$ws.Range("N86").Value = -3032
$ws.Range("H89").Value = 208.625
$ws.Range("I89").Value = 201.8
$ws.Range("J89").Value = 220
$ws.Range("K89").Value = 1816.2
$ws.Range("L89").Value = 1980
$ws.Range("M89").Value = 4111.8
$ws.Range("N89").Value = -13836
$ws.Range("H92").Value = 1030.3334
$ws.Range("I92").Value = 797
$ws.Range("J92").Value = 1497
$ws.Range("K92").Value = 2391
$ws.Range("L92").Value = 4491
$ws.Range("M92").Value = -1143
$ws.Range("N92").Value = -6987
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H113").Value = 1655.5
$ws.Range("I113").Value = 1655.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1655.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 514.5
$ws.Range("N113").ClearContents()
$ws.Range("H114").Value = 204361
$ws.Range("J114").Value = 204361
$ws.Range("L114").Value = 204361
$ws.Range("N114").Value = -213039
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 113
$ws.Range("I30").Value = 116.25
$ws.Range("J30").Value = 100
$ws.Range("K30").Value = 116.25
$ws.Range("L30").Value = 100
$ws.Range("M30").Value = -8.25
$ws.Range("N30").Value = -316
$ws.Range("H35").Value = 1101.8334
$ws.Range("I35").Value = 1222.2
$ws.Range("J35").Value = 500
$ws.Range("K35").Value = 1222.2
$ws.Range("L35").Value = 500
$ws.Range("M35").Value = -886.2
$ws.Range("N35").Value = -1172
$ws.Range("H39").Value = 23111.666
$ws.Range("I39").Value = 17559
$ws.Range("K39").Value = 17559
$ws.Range("M39").Value = -17099
$ws.Range("H128").Value = 39000
$ws.Range("J128").Value = 39000
$ws.Range("L128").Value = 39000
$ws.Range("N128").Value = -48960
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1342.4286
$ws.Range("I100").Value = 1174.25
$ws.Range("J100").Value = 1566.6666
$ws.Range("K100").Value = 2348.5
$ws.Range("L100").Value = 3133.3332
$ws.Range("M100").Value = -1807.5
$ws.Range("N100").Value = -4215.3332
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

Write-Host "Applied all changes"
